# Update timestamped test e-mail addresses in the "UsuariosRegistro" sheet
# from 20251111_202811 to 20251112_211458.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("UsuariosRegistro")

$ws.Range("C2").Value = "juan.perez+20251112_211458@test.com"
$ws.Range("C3").Value = "maria.gonzalez+20251112_211458@test.com"
$ws.Range("C4").Value = "carlos.rodriguez+20251112_211458@test.com"
$ws.Range("C5").Value = "ana.martinez+20251112_211458@test.com"
$ws.Range("C6").Value = "luis.garcia+20251112_211458@test.com"
